# Fixed size of images in some Excel solutions
#
# The document contains 6 inline pictures (in document order):
#   1. Picture 2  - screenshot of the Data Analysis button
#   2. Picture 1  - Random Number Generation menu item
#   3. Picture 10 - Random Number Generation dialog
#   4. Picture 11 - example generated data
#   5. Picture 12 - =RAND() column screenshot
#   6. Picture 13 - =NORM.INV() column screenshot
#
# Each one is being enlarged (roughly doubled) to make the screenshots
# easier to read. Sizes below are taken straight from the target EMU
# values (EMU / 12700 = points, the unit InlineShape.Width/Height use).

$d = $word.ActiveDocument
$shapes = $d.InlineShapes

$emuPerPoint = 12700.0

# id => [new width EMU, new height EMU]
$newSizesEmu = @{
    1 = @(1624273, 1243584)   # Picture 2
    2 = @(4572000, 2303630)   # Picture 1
    3 = @(4572000, 4172390)   # Picture 10
    4 = @(2119874, 4059936)   # Picture 11
    5 = @(3465647, 2743200)   # Picture 12
    6 = @(3850179, 2560320)   # Picture 13
}

for ($i = 1; $i -le $shapes.Count; $i++) {
    $shape = $shapes.Item($i)
    if ($newSizesEmu.ContainsKey($i)) {
        $size = $newSizesEmu[$i]
        $shape.Width = $size[0] / $emuPerPoint
        $shape.Height = $size[1] / $emuPerPoint
    }
}

Write-Output "Resized $($shapes.Count) inline pictures"
